$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data
$ws.Range("D2").Value = "'29.805.29"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "'1.620.99"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "'0.996"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'213.25"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'29.34"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  +8.94%  "
$ws.Range("D9").Value = "'0.258"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  +2.94%  "
$ws.Range("D10").Value = "'0.0606"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "'1.855.01"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "'1.627.09"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "'0.566"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  +5.39%  "
$ws.Range("D15").Value = "'3.91"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  +5.32%  "
$ws.Range("D16").Value = "'29.881.49"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "'8.79"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  +15.61%  "
$ws.Range("D18").Value = "'64.27"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'241.93"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'0.0₃0706"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("D21").Value = "'0.995"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").Value = "'9.57"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = "  +3.74%  "
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").Value = "'156.57"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "'15.59"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("D27").Value = "'0.110"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("D29").Value = "'0.997"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'0.0487"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  +3.21%  "
$ws.Range("D31").Value = "'1.11"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").Value = "'3.21"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "  +3.59%  "
$ws.Range("D34").Value = "'1.423.24"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +6.61%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").Value = "'2.86"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").Value = "'0.0169"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("D41").Value = "'0.0506"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "'1.97"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").Value = "'0.824"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("D44").Value = "'53.85"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "'69.03"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  +4.85%  "
$ws.Range("E46").Value = "  +19.05%  "
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Value = "'1.763.72"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  +1.18%  "

# Rows 50 and 51 swap ranking positions (BabyDogeCoin moves up to rank 50, Quant moves to rank 51)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₆0110"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  +6.81%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'88.01"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  +1.64%  "
